$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "26.991.32"
$ws.Range("E2").Value = "  +2.08%  "
$ws.Range("D3").Value = "1.816.74"
$ws.Range("E3").Value = "  +2.57%  "
Set-TextValue "D5" "311.68"
$ws.Range("E5").Value = "  +1.64%  "
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("E7").Value = "  -0.33%  "
Set-TextValue "D8" "0.3668"
$ws.Range("E8").Value = "  -0.04%  "
Set-TextValue "D9" "0.07261"
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("D10").Value = "2.167.43"
$ws.Range("E10").Value = "  +21.05%  "
Set-TextValue "D11" "0.8648"
$ws.Range("E11").Value = "  +1.57%  "
Set-TextValue "D12" "21.27"
$ws.Range("E12").Value = "  +4.76%  "
Set-TextValue "D13" "5.416"
$ws.Range("E13").Value = "  +3.28%  "
Set-TextValue "D14" "6.603"
$ws.Range("E14").Value = "  +2.54%  "
Set-TextValue "D15" "0.06948"
$ws.Range("E15").Value = "  +1.55%  "
Set-TextValue "D16" "81.13"
$ws.Range("E16").Value = "  +2.02%  "
Set-TextValue "D17" "1.013"
$ws.Range("E17").Value = "  +1.06%  "
Set-TextValue "D18" "0.000008802"
$ws.Range("E18").Value = "  +1.70%  "
$ws.Range("E19").Value = "  +0.41%  "
Set-TextValue "D20" "15.27"
$ws.Range("E20").Value = "  +1.41%  "
$ws.Range("D21").Value = "27.030.86"
$ws.Range("E21").Value = "  +2.23%  "
Set-TextValue "D22" "5.178"
$ws.Range("E22").Value = "  +1.03%  "
$ws.Range("D23").Value = "2.409.59"
$ws.Range("E23").Value = "  +19.79%  "
Set-TextValue "D24" "11.02"
$ws.Range("E24").Value = "  -2.09%  "
Set-TextValue "D25" "154.10"
$ws.Range("E25").Value = "  +1.36%  "
Set-TextValue "D26" "1.887"
$ws.Range("E26").Value = "  +1.79%  "
Set-TextValue "D27" "18.36"
$ws.Range("E27").Value = "  +1.17%  "
Set-TextValue "D28" "5.221"
$ws.Range("E28").Value = "  +2.42%  "
$ws.Range("E29").Value = "  +10.17%  "
Set-TextValue "D30" "114.69"
$ws.Range("E30").Value = "  -0.13%  "
Set-TextValue "D31" "0.08938"
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("E32").Value = "  +6.17%  "
Set-TextValue "D33" "0.7465"
$ws.Range("E33").Value = "  +2.98%  "
Set-TextValue "D34" "4.418"
$ws.Range("E34").Value = "  +2.05%  "
Set-TextValue "D35" "2.811"
$ws.Range("E35").Value = "  +2.32%  "
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("E37").Value = "  +4.71%  "
Set-TextValue "D38" "0.05215"
$ws.Range("E38").Value = "  +0.87%  "
Set-TextValue "D39" "0.01925"
$ws.Range("E39").Value = "  +1.58%  "
$ws.Range("E40").Value = "  +3.32%  "
Set-TextValue "D41" "0.1652"
$ws.Range("E41").Value = "  +2.92%  "
Set-TextValue "D42" "2.742"
$ws.Range("E42").Value = "  +8.22%  "
Set-TextValue "D43" "6.489"
$ws.Range("E43").Value = "  +3.99%  "
Set-TextValue "D44" "8.307"
$ws.Range("E44").Value = "  +3.35%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D45" "10.44"
$ws.Range("E45").Value = "  +2.24%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D46" "106.73"
$ws.Range("E46").Value = "  +1.69%  "
$ws.Range("E47").Value = "  +0.43%  "
Set-TextValue "D48" "0.4587"
$ws.Range("E48").Value = "  +2.36%  "
Set-TextValue "D49" "1.642"
$ws.Range("E49").Value = "  +3.56%  "
Set-TextValue "D50" "0.06218"
Set-TextValue "D51" "1.832"
$ws.Range("E51").Value = "  +4.79%  "
